$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.778.45'
$ws.Range('E2').Value = '  -0.68%  '
$ws.Range('D3').Value = '1.936.44'
$ws.Range('E3').Value = '  -0.88%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '243.03'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.99%  '
$ws.Range('E6').Value = '  -0.07%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4887'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.07%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2951'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.32%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06881'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.78%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.30'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.66%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '104.82'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -2.28%  '
$ws.Range('D12').Value = '1.938.21'
$ws.Range('E12').Value = '  -0.70%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.07779'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.44%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.338'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.26%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.7019'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.76%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '272.72'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -3.18%  '
$ws.Range('D17').Value = '30.788.45'
$ws.Range('E17').Value = '  -0.69%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '5.685'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +3.44%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000007718'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.41%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '13.11'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.88%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.000'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.09%  '
$ws.Range('B22').Value = 'BinanceUSD'
$ws.Range('C22').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.002'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.09%  '
$ws.Range('B23').Value = 'Chainlink'
$ws.Range('C23').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.540'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.82%  '
$ws.Range('B24').Value = 'Cosmos'
$ws.Range('C24').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '9.811'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.38%  '
$ws.Range('B25').Value = 'Monero'
$ws.Range('C25').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '164.78'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -2.88%  '
$ws.Range('B26').Value = 'EthereumClassic'
$ws.Range('C26').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '19.56'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -2.08%  '
$ws.Range('B27').Value = 'LidoDAOToken'
$ws.Range('C27').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.161'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.76%  '
$ws.Range('B28').Value = 'Stellar'
$ws.Range('C28').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.1035'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -2.02%  '
$ws.Range('B29').Value = 'Toncoin'
$ws.Range('C29').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.384'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -2.04%  '
$ws.Range('B30').Value = 'Filecoin'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.675'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +2.02%  '
$ws.Range('B31').Value = 'PancakeSwap'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.554'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.68%  '
$ws.Range('B32').Value = 'InternetComputer(DFINITY)'
$ws.Range('C32').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.420'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.86%  '
$ws.Range('B33').Value = 'Hedera'
$ws.Range('C33').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.04904'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.99%  '
$ws.Range('B34').Value = 'ImmutableX'
$ws.Range('C34').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.7584'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.84%  '
$ws.Range('B35').Value = 'ARBITRUM'
$ws.Range('C35').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.149'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.81%  '
$ws.Range('B36').Value = 'Frax'
$ws.Range('C36').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.000'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.07%  '
$ws.Range('B37').Value = 'HuobiToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.731'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.20%  '
$ws.Range('B38').Value = 'VeChain'
$ws.Range('C38').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.02004'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.28%  '
$ws.Range('B39').Value = 'MXToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.671'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -1.22%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '78.99'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +6.29%  '
$ws.Range('B41').Value = 'FraxShare'
$ws.Range('C41').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.454'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -1.46%  '
$ws.Range('B42').Value = 'RenderToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.076'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -3.31%  '
$ws.Range('B43').Value = 'TrustWalletToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.8978'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +1.60%  '
$ws.Range('B44').Value = 'TheSandbox'
$ws.Range('C44').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.4441'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.92%  '
$ws.Range('B45').Value = 'Quant'
$ws.Range('C45').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '108.26'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.29%  '
$ws.Range('B46').Value = 'Aptos'
$ws.Range('C46').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '7.869'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -3.36%  '
$ws.Range('B47').Value = 'PaxDollar'
$ws.Range('C47').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.000'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.09%  '
$ws.Range('B48').Value = 'Maker'
$ws.Range('C48').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '987.19'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.01%  '
$ws.Range('B49').Value = 'Algorand'
$ws.Range('C49').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.1248'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -1.24%  '
$ws.Range('B50').Value = 'Elrond'
$ws.Range('C50').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '36.21'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.25%  '
$ws.Range('B51').Value = 'EnergySwap'
$ws.Range('C51').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '9.244'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.63%  '
